# "Fruta / hortaliza, semanal" update:
# A new weekly price observation is inserted as a new row 72, pushing the
# existing rows 72-136 down to 73-137 (classic spreadsheet row insertion).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72; this shifts rows 72..136 down to 73..137
# and extends the used range from A1:T136 to A1:T137.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new observation.
$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "Vega Modelo de Temuco"
$ws.Range("C72").Value = "La Araucanía"
$ws.Range("D72").Value = 45240
$ws.Range("E72").Value = 9
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100108
$ws.Range("H72").Value = "Tropicales y subtropicales"
$ws.Range("I72").Value = 100108004
$ws.Range("J72").Value = "Papaya"
$ws.Range("K72").Value = "Cultivar IV Región"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 125
$ws.Range("N72").Value = 2600
$ws.Range("O72").Value = 2600
$ws.Range("P72").Value = 2600
$ws.Range("Q72").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R72").Value = "Provincia del Elquí"
$ws.Range("S72").Value = 2600
$ws.Range("T72").Value = 1
